# Final model run, saved results
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the discount_rate row values (row 2): lower/upper bounds 1.05 -> 5,
# unit changes from "absolute" to "percent".
$ws.Range("C2").Value = 5
$ws.Range("E2").Value = 5
$ws.Range("F2").Value = "percent"

# Update the saved selection to F3 (matches the workbook's last-saved cursor).
$ws.Range("F3").Select()
